# Apply the Alvearie FHIR IG "match-source-reference" StructureDefinition
# refresh (5.0.0 -> 6.0.0) to the workbook.
#
# Sheet 1 ("Metadata") is a Property/Value table:
#   - Version bumps 5.0.0 -> 6.0.0
#   - Date bumps to the new publish timestamp
#   - Publisher gets a real value ("Alvearie Team")
#   - The old duplicated "Contact" / "No display for ContactDetail" rows
#     are replaced by a single "Jurisdiction" / "United States of America"
#     row, and the leftover duplicate row is removed entirely (rows below
#     shift up by one).
#
# Sheet 2 ("Elements") is the per-element definition table; only the root
# Extension row's Short/Definition text changes to reflect the resource's
# own name/description instead of the generic placeholder text.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item(1)
$elements = $wb.Worksheets.Item(2)

# --- Sheet 1: Metadata -----------------------------------------------

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: new publish timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 becomes Jurisdiction / United States of America
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Remove the now-redundant duplicate "Contact" row (old row 11); this
# shifts every row below it up by one, matching the new A1:B20 dimension.
$meta.Rows.Item(11).Delete()

# --- Sheet 2: Elements -------------------------------------------------

# Root Extension element's Short/Definition reflect the real extension
# name & description now, instead of the generic placeholder text.
$elements.Range("K2").Value = "MatchSourceReference"
$elements.Range("L2").Value = "Reference to the source resource that is matched"
